# working_hours.xlsx - add a new pair of time entries (row 186/187) plus a
# trailing blank separator row, pushing the summary block (sum [min]/[h]/
# [working weeks]) down by two rows, and extend the F/G shared-formula
# fill range to cover the new data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make room: the old row 186 (blank separator) becomes real data, and two
# brand-new rows are inserted below it before the summary block.
$ws.Rows(187).Insert()
$ws.Rows(188).Insert()

# Row 186: new completed time entry (2014-08-18, 14:30-17:15)
$ws.Range("A186").Value = 2014
$ws.Range("B186").Value = 8
$ws.Range("C186").NumberFormat = "0"
$ws.Range("C186").Value = 18
$ws.Range("D186").Value = 0.60416666666666663
$ws.Range("E186").Value = 0.71875
$ws.Range("F186").Formula = "=(E186-D186)*24*60"
$ws.Range("G186").Formula = "=F186/60"

# Row 187: new, still-open entry (start time only, no end time yet)
$ws.Range("A187").Value = 2014
$ws.Range("B187").Value = 8
$ws.Range("C187").NumberFormat = "0"
$ws.Range("C187").Value = 18
$ws.Range("D187").Value = 0.84027777777777779

# Row 188 stays the blank separator row (formatting only, matches the
# pattern the old row 186 used before it got data).

# Summary block now lives on rows 189-191; point the totals at the
# extended data range (F2:F186 instead of F2:F185).
$ws.Range("F189").Formula = "=SUM(F2:F186)"
$ws.Range("F190").Formula = "=F189/60"
$ws.Range("F191").Formula = "=F190/38.5"

# Update the saved view state to match where the user left off editing.
[void]$ws.Range("E187").Select()

[void]$wb.Save()
